$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 425.0779544154803
$ws.Range("C2").Value = 0.0708109202757755
$ws.Range("D2").Value = 0.9580209895052474
$ws.Range("E2").Value = 0.9158759367194005
